$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text so values like "1.00" or "0.130" are not
# silently converted to numbers and lose their formatting, matching the original
# workbook where every Price cell is stored as an inline/text string.
$ws.Range("D1:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.005.05'
$ws.Range("E2").Value = '  -1.47%  '

$ws.Range("D3").Value = '3.563.94'
$ws.Range("E3").Value = '  -2.94%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '568.12'
$ws.Range("E5").Value = '  -5.19%  '

$ws.Range("D6").Value = '189.98'
$ws.Range("E6").Value = '  -1.08%  '

$ws.Range("D7").Value = '3.557.00'
$ws.Range("E7").Value = '  -3.07%  '

$ws.Range("D8").Value = '0.613'
$ws.Range("E8").Value = '  -1.52%  '

$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").Value = '0.674'
$ws.Range("E10").Value = '  -4.53%  '

$ws.Range("D11").Value = '55.67'
$ws.Range("E11").Value = '  -4.27%  '

$ws.Range("D12").Value = '0.148'
$ws.Range("E12").Value = '  -2.94%  '

$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -2.79%  '

$ws.Range("D14").Value = '9.82'
$ws.Range("E14").Value = '  -4.09%  '

$ws.Range("D15").Value = '4.140.08'
$ws.Range("E15").Value = '  -2.89%  '

$ws.Range("D16").Value = '3.573.43'
$ws.Range("E16").Value = '  -2.77%  '

$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("D18").Value = '66.960.45'
$ws.Range("E18").Value = '  -1.36%  '

$ws.Range("D19").Value = '12.14'
$ws.Range("E19").Value = '  -3.18%  '

$ws.Range("D20").Value = '18.12'
$ws.Range("E20").Value = '  -4.68%  '

$ws.Range("E21").Value = '  -5.64%  '

$ws.Range("D22").Value = '400.16'
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  -6.97%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '85.72'
$ws.Range("E24").Value = '  -2.70%  '

$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").Value = '11.83'
$ws.Range("E25").Value = '  +4.26%  '

$ws.Range("D26").Value = '2.90'
$ws.Range("E26").Value = '  -2.30%  '

$ws.Range("D27").Value = '12.43'
$ws.Range("E27").Value = '  -1.32%  '

$ws.Range("D28").Value = '6.09'
$ws.Range("E28").Value = '  +1.05%  '

$ws.Range("D29").Value = '3.64'
$ws.Range("E29").Value = '  -1.14%  '

$ws.Range("D30").Value = '7.78'
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("D31").Value = '8.93'
$ws.Range("E31").Value = '  -4.23%  '

$ws.Range("D32").Value = '31.08'
$ws.Range("E32").Value = '  -2.89%  '

$ws.Range("D33").Value = '637.23'
$ws.Range("E33").Value = '  +3.96%  '

$ws.Range("D34").Value = '12.08'
$ws.Range("E34").Value = '  -2.70%  '

$ws.Range("E35").Value = '  -3.51%  '

$ws.Range("D36").Value = '63.68'
$ws.Range("E36").Value = '  -6.32%  '

$ws.Range("D37").Value = '42.17'
$ws.Range("E37").Value = '  -7.50%  '

$ws.Range("D38").Value = '0.403'
$ws.Range("E38").Value = '  +1.28%  '

$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("D40").Value = '0.0₃0759'
$ws.Range("E40").Value = '  -2.28%  '

$ws.Range("D41").Value = '3.196.07'
$ws.Range("E41").Value = '  +12.93%  '

$ws.Range("D42").Value = '0.133'
$ws.Range("E42").Value = '  -1.74%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("E44").Value = '  +2.03%  '

$ws.Range("D45").Value = '2.69'
$ws.Range("E45").Value = '  +4.97%  '

$ws.Range("D46").Value = '0.0412'
$ws.Range("E46").Value = '  -3.68%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.14'
$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.130'
$ws.Range("E48").Value = '  -4.55%  '

$ws.Range("D49").Value = '141.51'
$ws.Range("E49").Value = '  -2.11%  '

$ws.Range("D50").Value = '8.47'
$ws.Range("E50").Value = '  -5.81%  '

$ws.Range("D51").Value = '2.52'
$ws.Range("E51").Value = '  -4.98%  '
